$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New package entry appended as row 4
$ws.Range("A4").Value = "GZIP"
$ws.Range("B4").Value = "GZIP"
$ws.Range("C4").Value = "GZIP"
$ws.Range("D4").Value = "1.0.0"
$ws.Range("E4").Value = "IFlow"

# F4 must hold the literal text "2025-08-08" (not an Excel date serial).
# Assigning the string directly gets auto-parsed into a date value, so
# write it as a formula returning the text, then convert that formula
# to its static value via copy / paste-special.
$ws.Range("F4").Formula = "=""2025-08-08"""
$ws.Range("F4").Copy()
$ws.Range("F4").PasteSpecial(-4163)
